$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2066.75
$ws.Range("I19").Value = 2391.625
$ws.Range("J19").Value = 1417
$ws.Range("K19").Value = 2391.625
$ws.Range("L19").Value = 1417
$ws.Range("M19").Value = -2216.625
$ws.Range("N19").Value = -1767
$ws.Range("H43").Value = 6943.3335
$ws.Range("H46").Value = 3000
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 9000
$ws.Range("N46").Value = -9238
$ws.Range("H60").Value = 3000
$ws.Range("J60").Value = 3000
$ws.Range("L60").Value = 9000
$ws.Range("N60").Value = -9968
$ws.Range("H98").Value = 613.0833
$ws.Range("I98").Value = 535.7
$ws.Range("K98").Value = 535.7
$ws.Range("M98").Value = 962.3
$ws.Range("H100").Value = 3335327.2
$ws.Range("I100").Value = 3335327.2
$ws.Range("K100").Value = 3335327.2
$ws.Range("M100").Value = -3334786.2
$ws.Range("H106").Value = 7999.3335
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H111").Value = 600
$ws.Range("I111").Value = 600
$ws.Range("K111").Value = 1800
$ws.Range("M111").Value = 1267
$ws.Range("H112").Value = 1781
$ws.Range("J112").Value = 1838.7059
$ws.Range("L112").Value = 5516.1177
$ws.Range("N112").Value = -7732.1177
$ws.Range("H116").Value = 9347.362999999999
$ws.Range("I116").Value = 8987.799999999999
$ws.Range("J116").Value = 9647
$ws.Range("K116").Value = 8987.799999999999
$ws.Range("L116").Value = 9647
$ws.Range("M116").Value = -5545.799999999999
$ws.Range("N116").Value = -16531
$ws.Range("H122").Value = 613.0833
$ws.Range("I122").Value = 535.7
$ws.Range("K122").Value = 1607.1
$ws.Range("M122").Value = 842.8999999999999
$ws.Range("H127").Value = 1374.75
$ws.Range("I127").Value = 500
$ws.Range("K127").Value = 1500
$ws.Range("M127").Value = 3460
$ws.Range("H137").Value = 2115.6191
$ws.Range("I137").Value = 1598.6364
$ws.Range("J137").Value = 2684.3
$ws.Range("K137").Value = 4795.9092
$ws.Range("L137").Value = 8052.900000000001
$ws.Range("M137").Value = -2245.9092
$ws.Range("N137").Value = -13152.9
$ws.Range("H138").Value = 3326.5894
$ws.Range("J138").Value = 3375.8542
$ws.Range("L138").Value = 10127.5626
$ws.Range("N138").Value = -20407.5626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1100.6666
$ws.Range("I2").Value = 1024.2
$ws.Range("K2").Value = 1024.2
$ws.Range("M2").Value = -911.2
$ws.Range("H32").Value = 8910.525
$ws.Range("I32").Value = 6138.4062
$ws.Range("K32").Value = 6138.4062
$ws.Range("M32").Value = -5851.4062
$ws.Range("H45").Value = 2137.1428
$ws.Range("I45").Value = 2192.2
$ws.Range("K45").Value = 2192.2
$ws.Range("M45").Value = -1815.2
$ws.Range("H74").Value = 799.25
$ws.Range("I74").Value = 732.3333
$ws.Range("K74").Value = 732.3333
$ws.Range("M74").Value = 141.6667
$ws.Range("H77").Value = 799.25
$ws.Range("I77").Value = 732.3333
$ws.Range("K77").Value = 3661.6665
$ws.Range("M77").Value = 706.3334999999997
$ws.Range("H116").Value = 1100.6666
$ws.Range("I116").Value = 1024.2
$ws.Range("K116").Value = 1024.2
$ws.Range("M116").Value = 1269.8
$ws.Range("H132").Value = 2967.3044
$ws.Range("I132").Value = 1912.3846
$ws.Range("K132").Value = 5737.1538
$ws.Range("M132").Value = -3207.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1100.6666
$ws.Range("I3").Value = 1024.2
$ws.Range("K3").Value = 1024.2
$ws.Range("M3").Value = -910.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1602.5
$ws.Range("I16").Value = 1602.5
$ws.Range("K16").Value = 1602.5
$ws.Range("M16").Value = -1315.5
$ws.Range("H58").Value = 2538.8235
$ws.Range("I58").Value = 2060.375
$ws.Range("J58").Value = 2964.111
$ws.Range("K58").Value = 2060.375
$ws.Range("L58").Value = 2964.111
$ws.Range("M58").Value = -1857.375
$ws.Range("N58").Value = -3370.111
$ws.Range("H99").Value = 4551.125
$ws.Range("I99").Value = 4399.3335
$ws.Range("K99").Value = 4399.3335
$ws.Range("M99").Value = -2901.3335
$ws.Range("H107").Value = 1870.5294
$ws.Range("I107").Value = 813.2
$ws.Range("K107").Value = 813.2
$ws.Range("M107").Value = 1106.8
$ws.Range("H113").Value = 1602.5
$ws.Range("I113").Value = 1602.5
$ws.Range("K113").Value = 1602.5
$ws.Range("M113").Value = 567.5
$ws.Range("H126").Value = 4551.125
$ws.Range("I126").Value = 4399.3335
$ws.Range("K126").Value = 13198.0005
$ws.Range("M126").Value = -10728.0005
$ws.Range("H134").Value = 3243.75
$ws.Range("I134").Value = 3243.75
$ws.Range("K134").Value = 9731.25
$ws.Range("M134").Value = -7196.25
$ws.Range("H136").Value = 2538.8235
$ws.Range("I136").Value = 2060.375
$ws.Range("J136").Value = 2964.111
$ws.Range("K136").Value = 6181.125
$ws.Range("L136").Value = 8892.332999999999
$ws.Range("M136").Value = -3631.125
$ws.Range("N136").Value = -13992.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 245.33333
$ws.Range("J107").Value = 245.33333
$ws.Range("L107").Value = 735.99999
$ws.Range("N107").Value = -4575.99999
$ws.Range("H132").Value = 3254
$ws.Range("J132").Value = 3778.8
$ws.Range("L132").Value = 34009.2
$ws.Range("N132").Value = -39069.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2284.8572
$ws.Range("I102").Value = 1999
$ws.Range("J102").Value = 2999.5
$ws.Range("K102").Value = 1999
$ws.Range("L102").Value = 2999.5
$ws.Range("M102").Value = -377
$ws.Range("N102").Value = -6243.5
$ws.Range("H113").Value = 1313.1
$ws.Range("I113").Value = 1313.1
$ws.Range("K113").Value = 1313.1
$ws.Range("M113").Value = 856.9000000000001
$ws.Range("H132").Value = 3017.5881
$ws.Range("I132").Value = 2256.2222
$ws.Range("K132").Value = 6768.6666
$ws.Range("M132").Value = -4238.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1001.3333
$ws.Range("I7").Value = 1001.3333
$ws.Range("K7").Value = 1001.3333
$ws.Range("M7").Value = -889.3333
$ws.Range("H126").Value = 1001.3333
$ws.Range("I126").Value = 1001.3333
$ws.Range("K126").Value = 3003.9999
$ws.Range("M126").Value = -533.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 619.9091
$ws.Range("I107").Value = 545.5714
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 1636.7142
$ws.Range("L107").Value = 2250
$ws.Range("M107").Value = 283.2857999999999
$ws.Range("N107").Value = -6090
$ws.Range("H122").Value = 2982.5833
$ws.Range("I122").Value = 3039.2
$ws.Range("K122").Value = 9117.599999999999
$ws.Range("M122").Value = -6667.599999999999
